$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recalculated "dSF" (column F) values for a handful of rows.
$ws.Range("F2").Value = -11
$ws.Range("F3").Value = -6
$ws.Range("F4").Value = 0
$ws.Range("F7").Value = -7
